$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Hoja1")

# Update "Salario Basico" for the first worker (row 16)
$ws.Range("G16").Value = 1600000

# Shift the "Periodo Mora" values for rows 17-23 so they read in
# ascending order (oldest -> newest) instead of descending, and
# move the "Valor Mora" amounts along with the periods they belong to.
$ws.Range("E17").Value = "2308"
$ws.Range("F17").Value = 46400

$ws.Range("E18").Value = "2309"
$ws.Range("F18").Value = 46400

$ws.Range("E19").Value = "2310"

$ws.Range("E20").Value = "2311"

$ws.Range("E21").Value = "2312"

$ws.Range("E22").Value = "2401"
$ws.Range("F22").Value = 40000

$ws.Range("E23").Value = "2402"
$ws.Range("F23").Value = 24000
